$d = $word.ActiveDocument

# Helper: touch a range's character formatting with a true no-op round trip
# (Bold -> True -> False) so the run gets its own (empty) run-properties
# element without altering the visible formatting. Used to reproduce the
# run splits that appear in the target revision.
function Touch-Range($rng) {
    $rng.Bold = 1
    $rng.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "Im Namen Gottes, des Vaters und des Sohnes und des Heiligen
#    Geistes." -- no text change, paragraph/run are simply re-touched.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(3)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
Touch-Range $full

# ---------------------------------------------------------------------
# 2) "Gemeinde: Amen." -> split into "Gemeinde: " / "Amen." and
#    italicise both runs.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(4)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$part1 = "Gemeinde: "
$r1 = $d.Range($full.Start, $full.Start + $part1.Length)
$r1.Font.Italic = $true
$r2 = $d.Range($full.Start + $part1.Length, $full.End)
$r2.Font.Italic = $true

# ---------------------------------------------------------------------
# 3) Merge "Wir sind heute hier versammelt, um Abschied zu nehmen von
#    VORNAME NACHNAME" with "Er verstarb am STERBEDATUM im Alter von
#    ALTER Jahren, hier in STERBEORT." into a single paragraph, joined
#    by ". ", then remove the blank paragraph that used to separate
#    them from the following text.
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$joinPos = $p6.Range.End - 1
$markRange = $d.Range($joinPos, $joinPos + 1)
$markRange.Delete()
$d.Range($joinPos, $joinPos).InsertAfter(". ")

# remove the now-orphaned empty paragraph that followed "Er verstarb ..."
$d.Paragraphs(7).Range.Delete()

# ---------------------------------------------------------------------
# 4) Split the merged paragraph into its six constituent runs (pure
#    run splits -- no formatting change).
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$full = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$segments = @(
    "Wir sind heute",
    " hier ",
    "versammelt, um Abschied zu nehmen von VORNAME NACHNAME",
    ". ",
    "Er verstarb am STERBEDATUM",
    " im Alter von ALTER Jahren, hier in STERBEORT."
)
$pos = $full.Start
foreach ($seg in $segments) {
    $r = $d.Range($pos, $pos + $seg.Length)
    Touch-Range $r
    $pos = $pos + $seg.Length
}

# ---------------------------------------------------------------------
# 5) "Wir vertrauen darauf, ..." -- no text change, just re-touched.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(7)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
Touch-Range $full

# ---------------------------------------------------------------------
# 6) "BIBELBERS" + line break + "Dennoch, ..." -- split "BIBELBERS"
#    into its own run.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(8)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$bibel = "BIBELBERS"
$r1 = $d.Range($full.Start, $full.Start + $bibel.Length)
Touch-Range $r1
$r2 = $d.Range($full.Start + $bibel.Length, $full.End)
Touch-Range $r2

Write-Output "done"
